$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.257.25"
$ws.Range("E2").Value = "  +2.65%  "
$ws.Range("D3").Value = "2.058.03"
$ws.Range("E3").Value = "  +1.55%  "
$ws.Range("E4").Value = "  +0.29%  "
$ws.Range("D5").Formula = "'228.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Formula = "'0.616"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.78%  "
$ws.Range("D7").Formula = "'60.96"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +8.71%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +2.08%  "
$ws.Range("D10").Formula = "'0.0827"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.72%  "
$ws.Range("E11").Value = "  +1.66%  "
$ws.Range("D12").Formula = "'14.80"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.50%  "
$ws.Range("D13").Value = "2.364.12"
$ws.Range("E13").Value = "  +1.67%  "
$ws.Range("D14").Formula = "'21.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.17%  "
$ws.Range("D15").Formula = "'0.761"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.90%  "
$ws.Range("E16").Value = "  +2.45%  "
$ws.Range("D17").Value = "2.054.19"
$ws.Range("E17").Value = "  +1.40%  "
$ws.Range("D18").Value = "38.178.57"
$ws.Range("E18").Value = "  +2.74%  "
$ws.Range("D19").Formula = "'6.23"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.06%  "
$ws.Range("D20").Formula = "'69.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.60%  "
$ws.Range("D21").Value = "0.0₃0835"
$ws.Range("E21").Value = "  +2.09%  "
$ws.Range("D22").Formula = "'225.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.04%  "
$ws.Range("D23").Formula = "'0.999"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("E25").Value = "  -0.88%  "
$ws.Range("D26").Formula = "'166.12"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.33%  "
$ws.Range("D27").Formula = "'9.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.79%  "
$ws.Range("E28").Value = "  +0.93%  "
$ws.Range("D29").Formula = "'18.97"
$ws.Range("D29").Style = "Normal"
$ws.Range("E30").Value = "  -1.73%  "
$ws.Range("E31").Value = "  +2.69%  "
$ws.Range("E32").Value = "  +0.52%  "
$ws.Range("B33").Value = "InternetComputer(DFINITY)"
$ws.Range("C33").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D33").Formula = "'4.56"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.96%  "
$ws.Range("B34").Value = "WEMIXToken"
$ws.Range("C34").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D34").Formula = "'2.05"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.61%  "
$ws.Range("D35").Formula = "'0.0604"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Formula = "'6.41"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +14.91%  "
$ws.Range("E37").Value = "  -1.88%  "
$ws.Range("D38").Formula = "'3.26"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.48%  "
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("D40").Value = "1.526.73"
$ws.Range("E40").Value = "  +3.91%  "
$ws.Range("D41").Formula = "'97.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +4.07%  "
$ws.Range("E42").Value = "  +1.72%  "
$ws.Range("D43").Formula = "'16.83"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.89%  "
$ws.Range("E44").Value = "  +1.55%  "
$ws.Range("E45").Value = "  +1.33%  "
$ws.Range("E46").Value = "  +2.14%  "
$ws.Range("E47").Value = "  -7.54%  "
$ws.Range("D48").Formula = "'1.01"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.58%  "
$ws.Range("E49").Value = "  +1.87%  "
$ws.Range("E50").Value = "  -0.47%  "
$ws.Range("D51").Value = "2.252.04"
$ws.Range("E51").Value = "  +1.85%  "
